$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: relabel column A, keep B/C text the same ---
$ws.Range("A13").Value = "13. Linked List"

# --- Row 14: relabel column A (and give it the same look as A13),
#     keep B the same, replace C's text and give C a new (non-wrapping) look ---
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "14. Linked List"

$ws.Range("B14").Value = "Remove duplicates"

$ws.Range("C14").Value = "Travers through the linked list. If curr.next == curr then skip over curr.next."
$ws.Range("C14").WrapText = $false
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").WrapText = $false
$ws.Range("C14").VerticalAlignment = -4160
$ws.Rows.Item(14).RowHeight = 22.5

$ws.Application.CutCopyMode = $false

# --- New row 15 ---
$ws.Range("A11").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "15. Linked list"

$ws.Range("B11").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "reverse a linked list"

$ws.Range("C12").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "Travers through the linked list. Save temp var for curr.next and update the curr.Next to the prev node and update curr and prev. At last return the prev"
$ws.Rows.Item(15).RowHeight = 26.25

# --- New row 16 ---
$ws.Range("A11").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "16. Linked list & two Ptrs"

$ws.Range("B11").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "merge two sorted lists"

$ws.Range("C16").Value = "go throw both lists and have one pointer for each list then check which node have the smallest value and add it to the new list, (increment the list pointer)"

# --- New row 17 ---
$ws.Range("A11").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "17. two pointers"

$ws.Range("B11").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "binary search"

$ws.Range("C17").Value = "use high low middle pointers if the middle is the target index then return it and if we don't find return -1 (low <= high while loop)"

# --- New row 18 ---
$ws.Range("A11").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "18. array"

$ws.Range("B11").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "smallest letter bigger then target"

$ws.Range("C18").Value = "return the first letter that bigger then the target if not found then return the first letter"

# --- New row 19 ---
$ws.Range("A11").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "19. Array & binary search"

$ws.Range("B11").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "find the peak of a mountain"

$ws.Range("C19").Value = "find the max of the array. Do the follow up."
$ws.Range("C19").Interior.Color = $ws.Range("C14").Interior.Color()

$ws.Application.CutCopyMode = $false

# --- Selection moves to the next empty row ---
$ws.Range("A20").Select()
